# Update the "within100" math-drill table cell contents.
# Values are given in row-major order (row 1 col 1..5, row 2 col 1..5, ...)
$newValues = @(
    "3+68=",
    "26+46=",
    "90-19=",
    "97-68=",
    "69+12=",
    "38-38=",
    "49-40=",
    "36+45=",
    "5+87=",
    "50-0=",
    "11+61=",
    "48+24=",
    "78-64=",
    "65+30=",
    "60-59=",
    "81-42=",
    "64-56=",
    "17+23=",
    "19-18=",
    "81+13=",
    "5-4=",
    "86+8=",
    "76-47=",
    "81-52=",
    "9+29=",
    "82-16=",
    "81-17=",
    "9+76=",
    "75-7=",
    "96-95=",
    "7+65=",
    "81-13=",
    "1+5=",
    "66+9=",
    "97-3=",
    "28+16=",
    "1+82=",
    "84-46=",
    "83-82=",
    "73-18=",
    "10-7=",
    "52-2=",
    "9+18=",
    "25+45=",
    "33+58=",
    "47+2=",
    "33+41=",
    "94-84=",
    "22-10=",
    "35-6=",
    "70+14=",
    "2+87=",
    "17-11=",
    "16+22=",
    "88+0=",
    "80-13=",
    "52-51=",
    "56-21=",
    "44-24=",
    "11+82=",
    "2+8=",
    "39+27=",
    "74-20=",
    "77+20=",
    "97+2=",
    "0+34=",
    "41+54=",
    "12-5=",
    "34+41=",
    "48+13=",
    "50+45=",
    "37+31=",
    "20+75=",
    "62-9=",
    "16+52=",
    "40-25=",
    "67-23=",
    "62-27=",
    "85-59=",
    "27+1=",
    "55+11=",
    "35-11=",
    "81+0=",
    "38+58=",
    "68+22=",
    "54-33=",
    "52-0=",
    "8+34=",
    "21+67=",
    "14+39=",
    "45+46=",
    "97-68=",
    "23+39=",
    "93-6=",
    "59-55=",
    "75+12=",
    "70-2=",
    "71-26=",
    "58-31=",
    "11+86="
)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$numRows = $t.Rows.Count
$numCols = $t.Columns.Count

$idx = 0
for ($r = 1; $r -le $numRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Updated $idx cells."
